$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Order matters: "11÷4=" is both an old value (-> "92÷3=") and a new value
# (target of "90÷4=" -> "11÷4="). Replace the pre-existing "11÷4=" cell
# first, before introducing a new "11÷4=" text elsewhere, so the second
# replacement does not re-match text we just wrote.
Replace-Exact "84÷3=" "83÷8="
Replace-Exact "23÷4=" "19÷7="
Replace-Exact "11÷4=" "92÷3="
Replace-Exact "52÷2=" "83÷6="
Replace-Exact "43÷6=" "17÷2="

Replace-Exact "22÷6=" "72÷6="
Replace-Exact "19÷2=" "93÷3="
Replace-Exact "78÷2=" "80÷4="
Replace-Exact "47÷3=" "81÷7="
Replace-Exact "63÷9=" "76÷5="

Replace-Exact "61÷9=" "47÷4="
Replace-Exact "90÷4=" "11÷4="
Replace-Exact "41÷3=" "88÷8="
Replace-Exact "23÷2=" "84÷5="
Replace-Exact "54÷6=" "16÷4="

Replace-Exact "95÷9=" "70÷9="
Replace-Exact "97÷9=" "48÷6="
Replace-Exact "76÷7=" "57÷9="
Replace-Exact "34÷7=" "17÷3="
Replace-Exact "74÷6=" "54÷2="

Replace-Exact "47÷8=" "96÷8="
Replace-Exact "14÷3=" "38÷4="
Replace-Exact "75÷8=" "55÷9="
Replace-Exact "99÷8=" "54÷9="
Replace-Exact "74÷4=" "89÷3="
